$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 280813.38
$ws.Range("J17").Value = 280813.38
$ws.Range("L17").Value = 842440.14
$ws.Range("N17").Value = -842776.14

$ws.Range("H43").Value = 2570
$ws.Range("I43").Value = 7950
$ws.Range("J43").Value = 1374.4445
$ws.Range("K43").Value = 7950
$ws.Range("L43").Value = 1374.4445
$ws.Range("M43").Value = -7881
$ws.Range("N43").Value = -1512.4445

$ws.Range("H138").Value = 3230.6667
$ws.Range("I138").Value = 1197.5938
$ws.Range("J138").Value = 4201.6865
$ws.Range("K138").Value = 3592.7814
$ws.Range("L138").Value = 12605.0595
$ws.Range("M138").Value = 1547.2186
$ws.Range("N138").Value = -22885.0595

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14845.281
$ws.Range("I32").Value = 8334.814
$ws.Range("J32").Value = 50001.8
$ws.Range("K32").Value = 8334.814
$ws.Range("L32").Value = 50001.8
$ws.Range("M32").Value = -8047.814
$ws.Range("N32").Value = -50575.8

$ws.Range("H61").Value = 1823.2609
$ws.Range("I61").Value = 1728.4412
$ws.Range("J61").Value = 2091.9167
$ws.Range("K61").Value = 1728.4412
$ws.Range("L61").Value = 2091.9167
$ws.Range("M61").Value = -1516.4412
$ws.Range("N61").Value = -2515.9167

$ws.Range("H102").Value = 5980.8335
$ws.Range("I102").Value = 3754.4443
$ws.Range("J102").Value = 12660
$ws.Range("K102").Value = 3754.4443
$ws.Range("L102").Value = 12660
$ws.Range("M102").Value = -2132.4443
$ws.Range("N102").Value = -15904

$ws.Range("H110").Value = 1504.6666
$ws.Range("I110").Value = 1500.5
$ws.Range("J110").Value = 1513
$ws.Range("K110").Value = 1500.5
$ws.Range("L110").Value = 1513
$ws.Range("M110").Value = 544.5
$ws.Range("N110").Value = -5603

$ws.Range("H132").Value = 1910.0698
$ws.Range("I132").Value = 1490.9
$ws.Range("J132").Value = 2877.3845
$ws.Range("K132").Value = 4472.700000000001
$ws.Range("L132").Value = 8632.1535
$ws.Range("M132").Value = -1942.700000000001
$ws.Range("N132").Value = -13692.1535

$ws.Range("H136").Value = 1823.2609
$ws.Range("I136").Value = 1728.4412
$ws.Range("J136").Value = 2091.9167
$ws.Range("K136").Value = 5185.3236
$ws.Range("L136").Value = 6275.750100000001
$ws.Range("M136").Value = -2635.3236
$ws.Range("N136").Value = -11375.7501

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H106").Value = 45000
$ws.Range("J106").Value = 45000
$ws.Range("L106").Value = 45000
$ws.Range("N106").Value = -47524

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1901.8334
$ws.Range("I16").Value = 1352.75
$ws.Range("J16").Value = 3000
$ws.Range("K16").Value = 1352.75
$ws.Range("L16").Value = 3000
$ws.Range("M16").Value = -1065.75
$ws.Range("N16").Value = -3574

$ws.Range("H99").Value = 1988.9778
$ws.Range("I99").Value = 1942.0646
$ws.Range("J99").Value = 2092.8572
$ws.Range("K99").Value = 1942.0646
$ws.Range("L99").Value = 2092.8572
$ws.Range("M99").Value = -444.0645999999999
$ws.Range("N99").Value = -5088.8572

$ws.Range("H113").Value = 1901.8334
$ws.Range("I113").Value = 1352.75
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 1352.75
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = 817.25
$ws.Range("N113").Value = -7340

$ws.Range("H126").Value = 1988.9778
$ws.Range("I126").Value = 1942.0646
$ws.Range("J126").Value = 2092.8572
$ws.Range("K126").Value = 5826.1938
$ws.Range("L126").Value = 6278.571599999999
$ws.Range("M126").Value = -3356.1938
$ws.Range("N126").Value = -11218.5716

$ws.Range("H132").Value = 1684.3158
$ws.Range("J132").Value = 2468
$ws.Range("L132").Value = 7404
$ws.Range("N132").Value = -12464

$ws.Range("H134").Value = 33334434
$ws.Range("I134").Value = 1206.5834
$ws.Range("K134").Value = 3619.7502
$ws.Range("M134").Value = -1084.7502

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H20").Value = 3300
$ws.Range("J20").Value = 3875
$ws.Range("L20").Value = 11625
$ws.Range("N20").Value = -12079

$ws.Range("H22").Value = 11111377
$ws.Range("J22").Value = 500
$ws.Range("L22").Value = 1500
$ws.Range("N22").Value = -1838

$ws.Range("H27").Value = 11111377
$ws.Range("J27").Value = 500
$ws.Range("L27").Value = 1500
$ws.Range("N27").Value = -1704

$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("N32").ClearContents()

$ws.Range("H39").Value = 4266.6665
$ws.Range("J39").Value = 4266.6665
$ws.Range("L39").Value = 12799.9995
$ws.Range("N39").Value = -13387.9995

$ws.Range("H49").Value = 4500
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 4500
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 13500
$ws.Range("M49").ClearContents()
$ws.Range("N49").Value = -13812

$ws.Range("H64").Value = 1790.2727
$ws.Range("I64").Value = 499.75
$ws.Range("J64").Value = 2527.7144
$ws.Range("K64").Value = 1499.25
$ws.Range("L64").Value = 7583.1432
$ws.Range("M64").Value = -1229.25
$ws.Range("N64").Value = -8123.1432

$ws.Range("H67").Value = 1790.2727
$ws.Range("I67").Value = 499.75
$ws.Range("J67").Value = 2527.7144
$ws.Range("K67").Value = 1499.25
$ws.Range("L67").Value = 7583.1432
$ws.Range("M67").Value = -563.25
$ws.Range("N67").Value = -9455.143199999999

$ws.Range("H109").Value = 18522040
$ws.Range("I109").Value = 62500590
$ws.Range("J109").Value = 4755.263
$ws.Range("K109").Value = 187501770
$ws.Range("L109").Value = 14265.789
$ws.Range("M109").Value = -187500730
$ws.Range("N109").Value = -16345.789

$ws.Range("H134").Value = 4965.2104
$ws.Range("I134").Value = 3987.6155
$ws.Range("J134").Value = 7083.3335
$ws.Range("K134").Value = 11962.8465
$ws.Range("L134").Value = 21250.0005
$ws.Range("M134").Value = -6892.8465
$ws.Range("N134").Value = -31390.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1435.8889
$ws.Range("I102").Value = 1365.375
$ws.Range("K102").Value = 1365.375
$ws.Range("M102").Value = 256.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1766.6666
$ws.Range("I46").Value = 2250
$ws.Range("J46").Value = 800
$ws.Range("K46").Value = 2250
$ws.Range("L46").Value = 800
$ws.Range("M46").Value = -2062
$ws.Range("N46").Value = -1176

$ws.Range("H132").Value = 3260.5454
$ws.Range("I132").Value = 2133.3333
$ws.Range("J132").Value = 4613.2
$ws.Range("K132").Value = 6399.999899999999
$ws.Range("L132").Value = 13839.6
$ws.Range("M132").Value = -3869.999899999999
$ws.Range("N132").Value = -18899.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H105").Value = 48000
$ws.Range("J105").Value = 48000
$ws.Range("L105").Value = 48000
$ws.Range("N105").Value = -54988

$ws.Range("H122").Value = 66668012
$ws.Range("I122").Value = 83334550
$ws.Range("J122").Value = 1862.6666
$ws.Range("K122").Value = 250003650
$ws.Range("L122").Value = 5587.9998
$ws.Range("M122").Value = -250001200
$ws.Range("N122").Value = -10487.9998

$ws.Range("H132").Value = 896.2273
$ws.Range("I132").Value = 715.3148
$ws.Range("J132").Value = 1710.3334
$ws.Range("K132").Value = 2145.9444
$ws.Range("L132").Value = 5131.0002
$ws.Range("M132").Value = 384.0556000000001
$ws.Range("N132").Value = -10191.0002
